$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 649
$ws.Range("F3").Value = 753
$ws.Range("F4").Value = 958
$ws.Range("F5").Value = 743
$ws.Range("F6").Value = 847
$ws.Range("F7").Value = 415
$ws.Range("F8").Value = 628
$ws.Range("F9").Value = 138
$ws.Range("F10").Value = 1230
$ws.Range("F11").Value = 658
$ws.Range("F14").Value = 170
$ws.Range("F15").Value = 18
$ws.Range("F16").Value = 634
$ws.Range("F18").Value = 369
$ws.Range("F19").Value = 360
$ws.Range("F20").Value = 84
$ws.Range("F21").Value = 560
$ws.Range("F22").Value = 99
$ws.Range("F23").Value = 592
$ws.Range("F25").Value = 821
$ws.Range("F26").Value = 8

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F8").Value = 185
$ws.Range("F9").Value = 228
$ws.Range("F11").Value = 27
$ws.Range("F12").Value = 24

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 649
$ws.Range("F7").Value = 753
$ws.Range("F8").Value = 958
$ws.Range("F9").Value = 743
$ws.Range("F10").Value = 847
$ws.Range("F11").Value = 415
$ws.Range("F12").Value = 628
$ws.Range("F13").Value = 138
$ws.Range("F14").Value = 1230
$ws.Range("F15").Value = 658
$ws.Range("F21").Value = 170
$ws.Range("F22").Value = 18
$ws.Range("F23").Value = 634
$ws.Range("F24").Value = 185
$ws.Range("F26").Value = 369
$ws.Range("F27").Value = 360
$ws.Range("F28").Value = 84
$ws.Range("F29").Value = 228
$ws.Range("F31").Value = 560
$ws.Range("F32").Value = 27
$ws.Range("F33").Value = 24
$ws.Range("F36").Value = 99
$ws.Range("F37").Value = 592
$ws.Range("F39").Value = 821
$ws.Range("F40").Value = 8
